$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 1731.3182
$ws.Range("I15").Value = 1731.3182
$ws.Range("K15").Value = 5193.9546
$ws.Range("M15").Value = -5024.9546

# Row 39
$ws.Range("H39").Value = 466.82352
$ws.Range("I39").Value = 178.66667
$ws.Range("K39").Value = 536.00001
$ws.Range("M39").Value = -240.00001

# Row 51
$ws.Range("H51").Value = 7669.2144
$ws.Range("I51").Value = 7560.846
$ws.Range("K51").Value = 7560.846
$ws.Range("M51").Value = -7076.846

# Row 95
$ws.Range("H95").Value = 51524.668
$ws.Range("J95").Value = 51524.668
$ws.Range("L95").Value = 51524.668
$ws.Range("N95").Value = -57016.668

# Row 116
$ws.Range("H116").Value = 18214.857
$ws.Range("I116").Value = 17999.5
$ws.Range("K116").Value = 17999.5
$ws.Range("M116").Value = -14557.5

# Row 129
$ws.Range("H129").Value = 1484.1904
$ws.Range("I129").Value = 878.4
$ws.Range("J129").Value = 2998.6667
$ws.Range("K129").Value = 2635.2
$ws.Range("L129").Value = 8996.000100000001
$ws.Range("M129").Value = 2364.8
$ws.Range("N129").Value = -18996.0001

# Row 135
$ws.Range("H135").Value = 961.26086
$ws.Range("I135").Value = 624.7143
$ws.Range("K135").Value = 5622.428699999999
$ws.Range("M135").Value = -3087.428699999999

# Row 138
$ws.Range("H138").Value = 2506.5974
$ws.Range("I138").Value = 1162.258
$ws.Range("J138").Value = 3412.5652
$ws.Range("K138").Value = 3486.774
$ws.Range("L138").Value = 10237.6956
$ws.Range("M138").Value = 1653.226
$ws.Range("N138").Value = -20517.6956


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7058.643
$ws.Range("I32").Value = 5858.25
$ws.Range("J32").Value = 10059.625
$ws.Range("K32").Value = 5858.25
$ws.Range("L32").Value = 10059.625
$ws.Range("M32").Value = -5571.25
$ws.Range("N32").Value = -10633.625

# Row 45
$ws.Range("H45").Value = 1620
$ws.Range("I45").Value = 1649.591
$ws.Range("K45").Value = 1649.591
$ws.Range("M45").Value = -1272.591

# Row 61
$ws.Range("H61").Value = 5503.4614
$ws.Range("I61").Value = 3811.0454
$ws.Range("J61").Value = 14811.75
$ws.Range("K61").Value = 3811.0454
$ws.Range("L61").Value = 14811.75
$ws.Range("M61").Value = -3599.0454
$ws.Range("N61").Value = -15235.75

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

# Row 135
$ws.Range("H135").Value = 60059.668
$ws.Range("J135").Value = 60059.668
$ws.Range("L135").Value = 60059.668
$ws.Range("N135").Value = -70199.668

# Row 136
$ws.Range("H136").Value = 5503.4614
$ws.Range("I136").Value = 3811.0454
$ws.Range("J136").Value = 14811.75
$ws.Range("K136").Value = 11433.1362
$ws.Range("L136").Value = 44435.25
$ws.Range("M136").Value = -8883.1362
$ws.Range("N136").Value = -49535.25


$ws = $wb.Worksheets.Item("BSM")
# Row 135
$ws.Range("H135").Value = 47207.31
$ws.Range("J135").Value = 47207.31
$ws.Range("L135").Value = 47207.31
$ws.Range("N135").Value = -57347.31


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 27708.6
$ws.Range("I31").Value = 2902.742
$ws.Range("J31").Value = 82635.86
$ws.Range("K31").Value = 2902.742
$ws.Range("L31").Value = 82635.86
$ws.Range("M31").Value = -2607.742
$ws.Range("N31").Value = -83225.86

# Row 34
$ws.Range("H34").Value = 27708.6
$ws.Range("I34").Value = 2902.742
$ws.Range("J34").Value = 82635.86
$ws.Range("K34").Value = 2902.742
$ws.Range("L34").Value = 82635.86
$ws.Range("M34").Value = -2700.742
$ws.Range("N34").Value = -83039.86

# Row 99
$ws.Range("H99").Value = 4975
$ws.Range("I99").Value = 4950
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 4950
$ws.Range("L99").Value = 5000
$ws.Range("M99").Value = -3452
$ws.Range("N99").Value = -7996

# Row 126
$ws.Range("H126").Value = 4975
$ws.Range("I126").Value = 4950
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 14850
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -12380
$ws.Range("N126").Value = -19940

# Row 134
$ws.Range("H134").Value = 3858.8572
$ws.Range("J134").Value = 9228
$ws.Range("L134").Value = 27684
$ws.Range("N134").Value = -32754

# Row 141
$ws.Range("H141").Value = 299743.75
$ws.Range("J141").Value = 299743.75
$ws.Range("L141").Value = 299743.75
$ws.Range("N141").Value = -310103.75


$ws = $wb.Worksheets.Item("CUL")
# Row 63
$ws.Range("H63").Value = 11662.667
$ws.Range("I63").Value = 4316.3335
$ws.Range("J63").Value = 19009
$ws.Range("K63").Value = 12949.0005
$ws.Range("L63").Value = 57027
$ws.Range("M63").Value = -12200.0005
$ws.Range("N63").Value = -58525

# Row 66
$ws.Range("H66").Value = 11662.667
$ws.Range("I66").Value = 4316.3335
$ws.Range("J66").Value = 19009
$ws.Range("K66").Value = 38847.0015
$ws.Range("L66").Value = 171081
$ws.Range("M66").Value = -35103.0015
$ws.Range("N66").Value = -178569

# Row 136
$ws.Range("H136").Value = 2664.25
$ws.Range("I136").Value = 2330.5715
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 6991.7145
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1891.7145
$ws.Range("N136").Value = -25200

# Row 141
$ws.Range("H141").Value = 10791.363
$ws.Range("I141").Value = 5647.25
$ws.Range("K141").Value = 16941.75
$ws.Range("M141").Value = -11761.75


$ws = $wb.Worksheets.Item("GSM")
# Row 108
$ws.Range("H108").Value = 75000
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 75000
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 75000
$ws.Range("M108").ClearContents()
$ws.Range("N108").Value = -82680

# Row 125
$ws.Range("H125").Value = 55000
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

# Row 126
$ws.Range("H126").Value = 4350.0586
$ws.Range("J126").Value = 5695.1
$ws.Range("L126").Value = 17085.3
$ws.Range("N126").Value = -22025.3

# Row 132
$ws.Range("H132").Value = 5755.237
$ws.Range("I132").Value = 5511.2666
$ws.Range("K132").Value = 16533.7998
$ws.Range("M132").Value = -14003.7998


$ws = $wb.Worksheets.Item("WVR")
# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 136
$ws.Range("H136").Value = 2361.06
$ws.Range("I136").Value = 2021.3877
$ws.Range("K136").Value = 6064.1631
$ws.Range("M136").Value = -3514.1631

# Row 137
$ws.Range("H137").Value = 69107.336
$ws.Range("J137").Value = 69107.336
$ws.Range("L137").Value = 69107.336
$ws.Range("N137").Value = -79307.336

